$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026462783288294
$ws.Range("D2").Value = 1.030756564931889
$ws.Range("E2").Value = 1.026685550906304
$ws.Range("F2").Value = 1.024979714447531
$ws.Range("I2").Value = 1.033926948399888
$ws.Range("J2").Value = 1.031626087283112
$ws.Range("K2").Value = 1.033566654395866
$ws.Range("L2").Value = 1.029507475098478
$ws.Range("M2").Value = 1.02780662802701
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027365618354493
$ws.Range("D3").Value = 1.031420536549801
$ws.Range("E3").Value = 1.027450649110859
$ws.Range("F3").Value = 1.026516456930325
$ws.Range("I3").Value = 1.034157352540494
$ws.Range("J3").Value = 1.032168932923749
$ws.Range("K3").Value = 1.03403943531443
$ws.Range("L3").Value = 1.030080254816403
$ws.Range("M3").Value = 1.029148595215377
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027949773309699
$ws.Range("D4").Value = 1.031850051738741
$ws.Range("E4").Value = 1.027946071304067
$ws.Range("F4").Value = 1.027510787652378
$ws.Range("I4").Value = 1.034305009878391
$ws.Range("J4").Value = 1.032519535329838
$ws.Range("K4").Value = 1.034344556731738
$ws.Range("L4").Value = 1.03045057565342
$ws.Range("M4").Value = 1.030016411982709
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028195342375544
$ws.Range("D5").Value = 1.032030590756867
$ws.Range("E5").Value = 1.028154430534757
$ws.Range("F5").Value = 1.027928798395322
$ws.Range("I5").Value = 1.034366742579107
$ws.Range("J5").Value = 1.032666771611646
$ws.Range("K5").Value = 1.034472638060112
$ws.Range("L5").Value = 1.030606184957996
$ws.Range("M5").Value = 1.030381120300874
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028236573929716
$ws.Range("D6").Value = 1.032060902321114
$ws.Range("E6").Value = 1.028189419858084
$ws.Range("F6").Value = 1.027998984030559
$ws.Range("I6").Value = 1.034377087684035
$ws.Range("J6").Value = 1.03269148402184
$ws.Range("K6").Value = 1.034494132226904
$ws.Range("L6").Value = 1.030632308120709
$ws.Range("M6").Value = 1.030442349403382
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027953054653821
$ws.Range("D7").Value = 1.031852464225458
$ws.Range("E7").Value = 1.027948855082087
$ws.Range("F7").Value = 1.027516373149636
$ws.Range("I7").Value = 1.034305836099595
$ws.Range("J7").Value = 1.032521503325231
$ws.Range("K7").Value = 1.034346268914446
$ws.Range("L7").Value = 1.030452655202092
$ws.Range("M7").Value = 1.030021285706911
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02676790874472
$ws.Range("D8").Value = 1.030980981428997
$ws.Range("E8").Value = 1.026944046053389
$ws.Range("F8").Value = 1.025499075819119
$ws.Range("I8").Value = 1.034005110431305
$ws.Range("J8").Value = 1.031809680043489
$ws.Range("K8").Value = 1.033726598258475
$ws.Range("L8").Value = 1.029701111834852
$ws.Range("M8").Value = 1.028260263395587
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024679233058536
$ws.Range("D9").Value = 1.029444438179077
$ws.Range("E9").Value = 1.025176171849323
$ws.Range("F9").Value = 1.021943787557476
$ws.Range("I9").Value = 1.033464253587397
$ws.Range("J9").Value = 1.030550346084387
$ws.Range("K9").Value = 1.032628550395407
$ws.Range("L9").Value = 1.028374467392577
$ws.Range("M9").Value = 1.025152910786019
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02328658907272
$ws.Range("D10").Value = 1.028419524776063
$ws.Range("E10").Value = 1.023999458685187
$ws.Range("F10").Value = 1.019572925543408
$ws.Range("I10").Value = 1.033096337683155
$ws.Range("J10").Value = 1.029707427672208
$ws.Range("K10").Value = 1.031892427943601
$ws.Range("L10").Value = 1.027488488944157
$ws.Range("M10").Value = 1.023078261464711
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022683513287154
$ws.Range("D11").Value = 1.02797560301884
$ws.Range("E11").Value = 1.023490379221881
$ws.Range("F11").Value = 1.01854609044973
$ws.Range("I11").Value = 1.032935285785256
$ws.Range("J11").Value = 1.029341637763392
$ws.Range("K11").Value = 1.031572711507771
$ws.Range("L11").Value = 1.027104486094066
$ws.Range("M11").Value = 1.022179127920689
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022459496226218
$ws.Range("D12").Value = 1.027810692120232
$ws.Range("E12").Value = 1.02330135179908
$ws.Range("F12").Value = 1.018164636868393
$ws.Range("I12").Value = 1.032875202322423
$ws.Range("J12").Value = 1.029205646771389
$ws.Range("K12").Value = 1.031453808932986
$ws.Range("L12").Value = 1.026961795162861
$ws.Range("M12").Value = 1.021845025030466
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022507549003476
$ws.Range("D13").Value = 1.027846066900202
$ws.Range("E13").Value = 1.023341895767205
$ws.Range("F13").Value = 1.018246461919438
$ws.Range("I13").Value = 1.032888102266611
$ws.Range("J13").Value = 1.02923482274685
$ws.Range("K13").Value = 1.031479320527596
$ws.Range("L13").Value = 1.026992405346208
$ws.Range("M13").Value = 1.021916696955489
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022664996127925
$ws.Range("D14").Value = 1.027961971808578
$ws.Range("E14").Value = 1.023474752779314
$ws.Range("F14").Value = 1.018514560254118
$ws.Range("I14").Value = 1.032930324603925
$ws.Range("J14").Value = 1.029330399160586
$ws.Range("K14").Value = 1.031562885954263
$ws.Range("L14").Value = 1.02709269234181
$ws.Range("M14").Value = 1.022151513443321
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022762003453185
$ws.Range("D15").Value = 1.028033382200956
$ws.Range("E15").Value = 1.023556619309088
$ws.Range("F15").Value = 1.018679738838269
$ws.Range("I15").Value = 1.032956304515033
$ws.Range("J15").Value = 1.029389270996883
$ws.Range("K15").Value = 1.031614354075604
$ws.Range("L15").Value = 1.027154475163248
$ws.Range("M15").Value = 1.022296174964197
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02332661212878
$ws.Range("D16").Value = 1.02844898376313
$ws.Range("E16").Value = 1.024033254031931
$ws.Range("F16").Value = 1.019641067733963
$ws.Range("I16").Value = 1.033106989474716
$ws.Range("J16").Value = 1.029731687071125
$ws.Range("K16").Value = 1.031913626035769
$ws.Range("L16").Value = 1.027513966202313
$ws.Range("M16").Value = 1.023137916788091
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023680762319576
$ws.Range("D17").Value = 1.028709645909577
$ws.Range("E17").Value = 1.024332353912603
$ws.Range("F17").Value = 1.020244016614213
$ws.Range("I17").Value = 1.033201043875502
$ws.Range("J17").Value = 1.029946261238073
$ws.Range("K17").Value = 1.032101091699155
$ws.Range("L17").Value = 1.027739366937071
$ws.Range("M17").Value = 1.023665702400573
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023887327147217
$ws.Range("D18").Value = 1.028861673391401
$ws.Range("E18").Value = 1.024506856726434
$ws.Range("F18").Value = 1.020595684097116
$ws.Range("I18").Value = 1.03325573614442
$ws.Range("J18").Value = 1.030071341469003
$ws.Range("K18").Value = 1.032210343593914
$ws.Range("L18").Value = 1.027870803822565
$ws.Range("M18").Value = 1.023973474283058
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023957759572034
$ws.Range("D19").Value = 1.028913508723978
$ws.Range("E19").Value = 1.024566364940434
$ws.Range("F19").Value = 1.020715589998599
$ws.Range("I19").Value = 1.033274356284343
$ws.Range("J19").Value = 1.03011397749354
$ws.Range("K19").Value = 1.032247579776228
$ws.Range("L19").Value = 1.027915614368989
$ws.Range("M19").Value = 1.024078403710352
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023642765872058
$ws.Range("D20").Value = 1.028681680600188
$ws.Range("E20").Value = 1.024300258883682
$ws.Range("F20").Value = 1.020179328263743
$ws.Range("I20").Value = 1.033190970112226
$ws.Range("J20").Value = 1.029923247451221
$ws.Range("K20").Value = 1.032080988093563
$ws.Range("L20").Value = 1.02771518724719
$ws.Range("M20").Value = 1.023609083951018
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022618632091806
$ws.Range("D21").Value = 1.027927841210687
$ws.Range("E21").Value = 1.023435627837393
$ws.Range("F21").Value = 1.018435613163892
$ws.Range("I21").Value = 1.032917898399225
$ws.Range("J21").Value = 1.02930225761106
$ws.Range("K21").Value = 1.031538282038463
$ws.Range("L21").Value = 1.027063161843012
$ws.Range("M21").Value = 1.022082369321428
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021974673000866
$ws.Range("D22").Value = 1.02745376516884
$ws.Range("E22").Value = 1.022892389920511
$ws.Range("F22").Value = 1.017339028380291
$ws.Range("I22").Value = 1.032744693777313
$ws.Range("J22").Value = 1.028911120416904
$ws.Range("K22").Value = 1.031196218560854
$ws.Range("L22").Value = 1.02665288844356
$ws.Range("M22").Value = 1.021121739534301
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022316052245996
$ws.Range("D23").Value = 1.027705091833812
$ws.Range("E23").Value = 1.023180333497116
$ws.Range("F23").Value = 1.017920373601281
$ws.Range("I23").Value = 1.032836656269762
$ws.Range("J23").Value = 1.029118535654377
$ws.Range("K23").Value = 1.031377632703189
$ws.Range("L23").Value = 1.026870412334082
$ws.Range("M23").Value = 1.021631057821887
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023659934846406
$ws.Range("D24").Value = 1.028694316957852
$ws.Range("E24").Value = 1.024314761113106
$ws.Range("F24").Value = 1.020208558212677
$ws.Range("I24").Value = 1.033195522531485
$ws.Range("J24").Value = 1.029933646629454
$ws.Range("K24").Value = 1.032090072336885
$ws.Range("L24").Value = 1.027726113118009
$ws.Range("M24").Value = 1.023634667626548
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025219240338076
$ws.Range("D25").Value = 1.02984177148902
$ws.Range("E25").Value = 1.025632882754793
$ws.Range("F25").Value = 1.022863010937346
$ws.Range("I25").Value = 1.033605373362013
$ws.Range("J25").Value = 1.030876507187855
$ws.Range("K25").Value = 1.032913144282034
$ws.Range("L25").Value = 1.028717710997193
$ws.Range("M25").Value = 1.025956760369553
